# Matchmaking integration (administrator interface)
# Adds the "MC" (match list) request/response protocol rows to the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")
$ws.Activate()

# NOTE: new shared-string entries are written in the order they must end up
# in sharedStrings.xml (uniqueCount indices 46, 47, 48):
#   46 "MC"
#   47 "Envoi la liste des matchs planifiés et en cours (interface admin)"
#   48 "Demande de la liste des matchs créés (planifiés et en cours) (interface admin)"
$ws.Range("C20").Value = "MC"
$ws.Range("E21").Value = "Envoi la liste des matchs planifiés et en cours (interface admin)"
$ws.Range("E20").Value = "Demande de la liste des matchs créés (planifiés et en cours) (interface admin)"

# Row 20: Client -> Serveur : MC request (list of created matches)
$ws.Range("A20").Value = "Client"
$ws.Range("B20").Value = "Serveur"

# Row 21: Serveur -> Client : MC response (list of planned/ongoing matches)
$ws.Range("A21").Value = "Serveur"
$ws.Range("B21").Value = "Client"
$ws.Range("C21").Value = "MC"
$ws.Range("D21").Value = "(voir code)"

# Both new rows wrap to two lines of text, same as the other multi-line rows.
$ws.Rows("20:21").RowHeight = 30

# Move the selection down to where the user ended up after the edit.
$ws.Range("A22").Select() | Out-Null
